# Generate Report for Handback
# Update the timestamps recorded for the d08ad933-095b-40df-b7c4-daa762fa9ddc
# handback entry across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-18 12:46:06"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-18 12:45:55"
$wsZhCn.Range("K4").Value = "2016-08-18 12:46:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-18 12:46:36"
